$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 12-27 and add new rows 28-34 per the diff
# Ensure new date cells (rows 28-34, column D) use the same date format as existing date cells
$ws.Range("D28:D34").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 12
$ws.Range("A12").Value = 6
$ws.Range("B12").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C12").Value = "Metropolitana"
$ws.Range("D12").Value = 44482
$ws.Range("E12").Value = 13
$ws.Range("F12").Value = 300000000
$ws.Range("G12").Value = "Espárragos"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Banquete"
$ws.Range("J12").Value = 580
$ws.Range("K12").Value = 1500
$ws.Range("L12").Value = 1600
$ws.Range("M12").Value = 1560
$ws.Range("N12").Value = "`$/kilo"
$ws.Range("O12").Value = "Provincia de Linares"
$ws.Range("P12").Value = 1560
$ws.Range("Q12").Value = 1
$ws.Range("R12").Value = "Hortaliza"

# Row 13
$ws.Range("A13").Value = 6
$ws.Range("B13").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C13").Value = "Metropolitana"
$ws.Range("D13").Value = 44482
$ws.Range("E13").Value = 13
$ws.Range("F13").Value = 300000000
$ws.Range("G13").Value = "Espárragos"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Banquete"
$ws.Range("J13").Value = 920
$ws.Range("K13").Value = 1500
$ws.Range("L13").Value = 1600
$ws.Range("M13").Value = 1565
$ws.Range("N13").Value = "`$/kilo"
$ws.Range("O13").Value = "Región Metropolitana"
$ws.Range("P13").Value = 1565
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = "Hortaliza"

# Row 14
$ws.Range("A14").Value = 6
$ws.Range("B14").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C14").Value = "Metropolitana"
$ws.Range("D14").Value = 44482
$ws.Range("E14").Value = 13
$ws.Range("F14").Value = 300000000
$ws.Range("G14").Value = "Espárragos"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 520
$ws.Range("K14").Value = 1300
$ws.Range("L14").Value = 1400
$ws.Range("M14").Value = 1348
$ws.Range("N14").Value = "`$/kilo"
$ws.Range("O14").Value = "Provincia de Linares"
$ws.Range("P14").Value = 1348
$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = "Hortaliza"

# Row 15
$ws.Range("A15").Value = 6
$ws.Range("B15").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C15").Value = "Metropolitana"
$ws.Range("D15").Value = 44482
$ws.Range("E15").Value = 13
$ws.Range("F15").Value = 300000000
$ws.Range("G15").Value = "Espárragos"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 660
$ws.Range("K15").Value = 1300
$ws.Range("L15").Value = 1400
$ws.Range("M15").Value = 1361
$ws.Range("N15").Value = "`$/kilo"
$ws.Range("O15").Value = "Región Metropolitana"
$ws.Range("P15").Value = 1361
$ws.Range("Q15").Value = 1
$ws.Range("R15").Value = "Hortaliza"

# Row 16
$ws.Range("A16").Value = 6
$ws.Range("B16").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C16").Value = "Metropolitana"
$ws.Range("D16").Value = 44482
$ws.Range("E16").Value = 13
$ws.Range("F16").Value = 300000000
$ws.Range("G16").Value = "Espárragos"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Segunda"
$ws.Range("J16").Value = 350
$ws.Range("K16").Value = 1100
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = 1143
$ws.Range("N16").Value = "`$/kilo"
$ws.Range("O16").Value = "Provincia de Linares"
$ws.Range("P16").Value = 1143
$ws.Range("Q16").Value = 1
$ws.Range("R16").Value = "Hortaliza"

# Row 17
$ws.Range("A17").Value = 6
$ws.Range("B17").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C17").Value = "Metropolitana"
$ws.Range("D17").Value = 44482
$ws.Range("E17").Value = 13
$ws.Range("F17").Value = 300000000
$ws.Range("G17").Value = "Espárragos"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Segunda"
$ws.Range("J17").Value = 470
$ws.Range("K17").Value = 1100
$ws.Range("L17").Value = 1200
$ws.Range("M17").Value = 1164
$ws.Range("N17").Value = "`$/kilo"
$ws.Range("O17").Value = "Región Metropolitana"
$ws.Range("P17").Value = 1164
$ws.Range("Q17").Value = 1
$ws.Range("R17").Value = "Hortaliza"

# Row 18
$ws.Range("A18").Value = 6
$ws.Range("B18").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C18").Value = "Metropolitana"
$ws.Range("D18").Value = 44482
$ws.Range("E18").Value = 13
$ws.Range("F18").Value = 300000000
$ws.Range("G18").Value = "Espárragos"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Tercera"
$ws.Range("J18").Value = 550
$ws.Range("K18").Value = 800
$ws.Range("L18").Value = 900
$ws.Range("M18").Value = 842
$ws.Range("N18").Value = "`$/kilo"
$ws.Range("O18").Value = "Provincia de Linares"
$ws.Range("P18").Value = 842
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = "Hortaliza"

# Row 19
$ws.Range("A19").Value = 6
$ws.Range("B19").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C19").Value = "Metropolitana"
$ws.Range("D19").Value = 44467
$ws.Range("E19").Value = 13
$ws.Range("F19").Value = 300000000
$ws.Range("G19").Value = "Espárragos"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Banquete"
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 2000
$ws.Range("L19").Value = 2000
$ws.Range("M19").Value = 2000
$ws.Range("N19").Value = "`$/caja 10 kilos"
$ws.Range("O19").Value = "Provincia de Linares"
$ws.Range("P19").Value = 200
$ws.Range("Q19").Value = 10
$ws.Range("R19").Value = "Hortaliza"

# Row 20
$ws.Range("A20").Value = 6
$ws.Range("B20").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C20").Value = "Metropolitana"
$ws.Range("D20").Value = 44467
$ws.Range("E20").Value = 13
$ws.Range("F20").Value = 300000000
$ws.Range("G20").Value = "Espárragos"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 150
$ws.Range("K20").Value = 1500
$ws.Range("L20").Value = 1500
$ws.Range("M20").Value = 1500
$ws.Range("N20").Value = "`$/caja 10 kilos"
$ws.Range("O20").Value = "Provincia de Linares"
$ws.Range("P20").Value = 150
$ws.Range("Q20").Value = 10
$ws.Range("R20").Value = "Hortaliza"

# Row 21
$ws.Range("A21").Value = 6
$ws.Range("B21").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C21").Value = "Metropolitana"
$ws.Range("D21").Value = 44467
$ws.Range("E21").Value = 13
$ws.Range("F21").Value = 300000000
$ws.Range("G21").Value = "Espárragos"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Segunda"
$ws.Range("J21").Value = 50
$ws.Range("K21").Value = 1300
$ws.Range("L21").Value = 1300
$ws.Range("M21").Value = 1300
$ws.Range("N21").Value = "`$/caja 10 kilos"
$ws.Range("O21").Value = "Provincia de Linares"
$ws.Range("P21").Value = 130
$ws.Range("Q21").Value = 10
$ws.Range("R21").Value = "Hortaliza"

# Row 22
$ws.Range("A22").Value = 6
$ws.Range("B22").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C22").Value = "Metropolitana"
$ws.Range("D22").Value = 44161
$ws.Range("E22").Value = 13
$ws.Range("F22").Value = 300000000
$ws.Range("G22").Value = "Espárragos"
$ws.Range("H22").Value = "Verde"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 4300
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = 1000
$ws.Range("N22").Value = "`$/kilo"
$ws.Range("O22").Value = "Provincia de Linares"
$ws.Range("P22").Value = 1000
$ws.Range("Q22").Value = 1
$ws.Range("R22").Value = "Hortaliza"

# Row 23
$ws.Range("A23").Value = 6
$ws.Range("B23").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C23").Value = "Metropolitana"
$ws.Range("D23").Value = 44161
$ws.Range("E23").Value = 13
$ws.Range("F23").Value = 300000000
$ws.Range("G23").Value = "Espárragos"
$ws.Range("H23").Value = "Verde"
$ws.Range("I23").Value = "Segunda"
$ws.Range("J23").Value = 2500
$ws.Range("K23").Value = 800
$ws.Range("L23").Value = 800
$ws.Range("M23").Value = 800
$ws.Range("N23").Value = "`$/kilo"
$ws.Range("O23").Value = "Provincia de Linares"
$ws.Range("P23").Value = 800
$ws.Range("Q23").Value = 1
$ws.Range("R23").Value = "Hortaliza"

# Row 24
$ws.Range("A24").Value = 6
$ws.Range("B24").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C24").Value = "Metropolitana"
$ws.Range("D24").Value = 44474
$ws.Range("E24").Value = 13
$ws.Range("F24").Value = 300000000
$ws.Range("G24").Value = "Espárragos"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Banquete"
$ws.Range("J24").Value = 780
$ws.Range("K24").Value = 1500
$ws.Range("L24").Value = 1600
$ws.Range("M24").Value = 1558
$ws.Range("N24").Value = "`$/kilo"
$ws.Range("O24").Value = "Provincia de Linares"
$ws.Range("P24").Value = 1558
$ws.Range("Q24").Value = 1
$ws.Range("R24").Value = "Hortaliza"

# Row 25
$ws.Range("A25").Value = 6
$ws.Range("B25").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C25").Value = "Metropolitana"
$ws.Range("D25").Value = 44474
$ws.Range("E25").Value = 13
$ws.Range("F25").Value = 300000000
$ws.Range("G25").Value = "Espárragos"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 520
$ws.Range("K25").Value = 1300
$ws.Range("L25").Value = 1400
$ws.Range("M25").Value = 1348
$ws.Range("N25").Value = "`$/kilo"
$ws.Range("O25").Value = "Provincia de Linares"
$ws.Range("P25").Value = 1348
$ws.Range("Q25").Value = 1
$ws.Range("R25").Value = "Hortaliza"

# Row 26
$ws.Range("A26").Value = 6
$ws.Range("B26").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C26").Value = "Metropolitana"
$ws.Range("D26").Value = 44474
$ws.Range("E26").Value = 13
$ws.Range("F26").Value = 300000000
$ws.Range("G26").Value = "Espárragos"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Segunda"
$ws.Range("J26").Value = 400
$ws.Range("K26").Value = 1000
$ws.Range("L26").Value = 1200
$ws.Range("M26").Value = 1100
$ws.Range("N26").Value = "`$/kilo"
$ws.Range("O26").Value = "Provincia de Linares"
$ws.Range("P26").Value = 1100
$ws.Range("Q26").Value = 1
$ws.Range("R26").Value = "Hortaliza"

# Row 27
$ws.Range("A27").Value = 6
$ws.Range("B27").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C27").Value = "Metropolitana"
$ws.Range("D27").Value = 44159
$ws.Range("E27").Value = 13
$ws.Range("F27").Value = 300000000
$ws.Range("G27").Value = "Espárragos"
$ws.Range("H27").Value = "Verde"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 4300
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = 1000
$ws.Range("N27").Value = "`$/kilo"
$ws.Range("O27").Value = "Región Metropolitana"
$ws.Range("P27").Value = 1000
$ws.Range("Q27").Value = 1
$ws.Range("R27").Value = "Hortaliza"

# Row 28
$ws.Range("A28").Value = 6
$ws.Range("B28").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C28").Value = "Metropolitana"
$ws.Range("D28").Value = 44159
$ws.Range("E28").Value = 13
$ws.Range("F28").Value = 300000000
$ws.Range("G28").Value = "Espárragos"
$ws.Range("H28").Value = "Verde"
$ws.Range("I28").Value = "Segunda"
$ws.Range("J28").Value = 2500
$ws.Range("K28").Value = 800
$ws.Range("L28").Value = 800
$ws.Range("M28").Value = 800
$ws.Range("N28").Value = "`$/kilo"
$ws.Range("O28").Value = "Región Metropolitana"
$ws.Range("P28").Value = 800
$ws.Range("Q28").Value = 1
$ws.Range("R28").Value = "Hortaliza"

# Row 29
$ws.Range("A29").Value = 6
$ws.Range("B29").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C29").Value = "Metropolitana"
$ws.Range("D29").Value = 44476
$ws.Range("E29").Value = 13
$ws.Range("F29").Value = 300000000
$ws.Range("G29").Value = "Espárragos"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Banquete"
$ws.Range("J29").Value = 1300
$ws.Range("K29").Value = 1500
$ws.Range("L29").Value = 1600
$ws.Range("M29").Value = 1554
$ws.Range("N29").Value = "`$/kilo"
$ws.Range("O29").Value = "Provincia de Linares"
$ws.Range("P29").Value = 1554
$ws.Range("Q29").Value = 1
$ws.Range("R29").Value = "Hortaliza"

# Row 30
$ws.Range("A30").Value = 6
$ws.Range("B30").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C30").Value = "Metropolitana"
$ws.Range("D30").Value = 44476
$ws.Range("E30").Value = 13
$ws.Range("F30").Value = 300000000
$ws.Range("G30").Value = "Espárragos"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Banquete"
$ws.Range("J30").Value = 700
$ws.Range("K30").Value = 1400
$ws.Range("L30").Value = 1500
$ws.Range("M30").Value = 1457
$ws.Range("N30").Value = "`$/kilo"
$ws.Range("O30").Value = "Región Metropolitana"
$ws.Range("P30").Value = 1457
$ws.Range("Q30").Value = 1
$ws.Range("R30").Value = "Hortaliza"

# Row 31
$ws.Range("A31").Value = 6
$ws.Range("B31").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C31").Value = "Metropolitana"
$ws.Range("D31").Value = 44476
$ws.Range("E31").Value = 13
$ws.Range("F31").Value = 300000000
$ws.Range("G31").Value = "Espárragos"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 900
$ws.Range("K31").Value = 1300
$ws.Range("L31").Value = 1400
$ws.Range("M31").Value = 1356
$ws.Range("N31").Value = "`$/kilo"
$ws.Range("O31").Value = "Provincia de Linares"
$ws.Range("P31").Value = 1356
$ws.Range("Q31").Value = 1
$ws.Range("R31").Value = "Hortaliza"

# Row 32
$ws.Range("A32").Value = 6
$ws.Range("B32").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C32").Value = "Metropolitana"
$ws.Range("D32").Value = 44476
$ws.Range("E32").Value = 13
$ws.Range("F32").Value = 300000000
$ws.Range("G32").Value = "Espárragos"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 500
$ws.Range("K32").Value = 1200
$ws.Range("L32").Value = 1300
$ws.Range("M32").Value = 1260
$ws.Range("N32").Value = "`$/kilo"
$ws.Range("O32").Value = "Región Metropolitana"
$ws.Range("P32").Value = 1260
$ws.Range("Q32").Value = 1
$ws.Range("R32").Value = "Hortaliza"

# Row 33
$ws.Range("A33").Value = 6
$ws.Range("B33").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C33").Value = "Metropolitana"
$ws.Range("D33").Value = 44476
$ws.Range("E33").Value = 13
$ws.Range("F33").Value = 300000000
$ws.Range("G33").Value = "Espárragos"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Segunda"
$ws.Range("J33").Value = 500
$ws.Range("K33").Value = 1100
$ws.Range("L33").Value = 1200
$ws.Range("M33").Value = 1160
$ws.Range("N33").Value = "`$/kilo"
$ws.Range("O33").Value = "Provincia de Linares"
$ws.Range("P33").Value = 1160
$ws.Range("Q33").Value = 1
$ws.Range("R33").Value = "Hortaliza"

# Row 34
$ws.Range("A34").Value = 6
$ws.Range("B34").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C34").Value = "Metropolitana"
$ws.Range("D34").Value = 44476
$ws.Range("E34").Value = 13
$ws.Range("F34").Value = 300000000
$ws.Range("G34").Value = "Espárragos"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Segunda"
$ws.Range("J34").Value = 200
$ws.Range("K34").Value = 1000
$ws.Range("L34").Value = 1100
$ws.Range("M34").Value = 1050
$ws.Range("N34").Value = "`$/kilo"
$ws.Range("O34").Value = "Región Metropolitana"
$ws.Range("P34").Value = 1050
$ws.Range("Q34").Value = 1
$ws.Range("R34").Value = "Hortaliza"

